$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B3").Value = "wlthrcde_fctb"
$ws.Range("B3").Select()
